# Append new TC / Customer_ID / PD rows (rows 273-289) to the "Individual
# Customers" data sheet, per the R22 UAT2 (Regression) test-data commit.
#
# The source values look numeric (e.g. "17866751"), but the existing sheet
# stores every cell in columns A:C as a shared string (t="s"), never as a
# native number. Assigning them directly via .Value/.Value2 would make Excel
# auto-coerce the digit-only strings into numbers, so instead each value is
# written as a literal-text formula (="17866751") and then the whole new
# block is copied and pasted back as values. PasteSpecial(values) keeps the
# already-evaluated text results (turning them into plain shared strings)
# without leaving behind a formula or any new cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("118500", "17866751", "6004"),
    @("118518", "17866752", "6020"),
    @("118498", "17866753", "1001"),
    @("118452", "17866754", "1001"),
    @("118518", "17866755", "1035"),
    @("118452", "17866756", "1150"),
    @("118448", "17866757", "1068"),
    @("118448", "17866759", "1005"),
    @("118500", "17866760", "6004"),
    @("118518", "17866761", "6020"),
    @("118498", "17866762", "1001"),
    @("118452", "17866763", "1001"),
    @("118518", "17866764", "6005"),
    @("118448", "17866766", "1047"),
    @("118452", "17866767", "1150"),
    @("118448", "17866768", "1068"),
    @("118448", "17866770", "1005")
)

$startRow = 273
$endRow = $startRow + $newRows.Length - 1

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).Formula = "=""" + $vals[0] + """"
    $ws.Cells.Item($r, 2).Formula = "=""" + $vals[1] + """"
    $ws.Cells.Item($r, 3).Formula = "=""" + $vals[2] + """"
}

$fillRange = $ws.Range("A" + $startRow + ":C" + $endRow)
$fillRange.Copy()
$fillRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0
